$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.104248642921448
$ws.Range("B1").Value = 2.453152656555176
$ws.Range("C1").Value = 5.276597023010254
$ws.Range("D1").Value = 2.217790126800537
$ws.Range("E1").Value = 1.277986168861389
